$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New working set of verb stimuli for column B ("word"), rows 2-193.
$words = @(
    "steuern",
    "klappen",
    "knarren",
    "schalten",
    "holen",
    "scheitern",
    "schaden",
    "mögen",
    "warnen",
    "trösten",
    "melden",
    "rasen",
    "boxen",
    "wirken",
    "duschen",
    "zeigen",
    "erben",
    "mauern",
    "runden",
    "kehren",
    "tropfen",
    "wellen",
    "flüstern",
    "zielen",
    "pressen",
    "opfern",
    "sorgen",
    "räumen",
    "stellen",
    "tauschen",
    "schätzen",
    "liefern",
    "arten",
    "doppeln",
    "gründen",
    "lockern",
    "hauen",
    "jagen",
    "kichern",
    "proben",
    "freuen",
    "pflanzen",
    "tollen",
    "fischen",
    "schultern",
    "treiben",
    "weichen",
    "schreiten",
    "heben",
    "drehen",
    "schicken",
    "fordern",
    "öffnen",
    "folgen",
    "starren",
    "münzen",
    "schneiden",
    "ärgern",
    "fallen",
    "stecken",
    "sperren",
    "dringen",
    "lohnen",
    "sinken",
    "spannen",
    "feiern",
    "gelten",
    "wehtun",
    "flehen",
    "tauchen",
    "fließen",
    "suchen",
    "morden",
    "wahren",
    "flüchten",
    "malen",
    "reizen",
    "bauen",
    "drücken",
    "lesen",
    "dauern",
    "fällen",
    "bluten",
    "streichen",
    "dienen",
    "altern",
    "greifen",
    "fassen",
    "wandern",
    "seufzen",
    "mühen",
    "testen",
    "lehnen",
    "scheinen",
    "stechen",
    "enden",
    "geben",
    "platzen",
    "klingeln",
    "trauen",
    "kosten",
    "filmen",
    "schwimmen",
    "töten",
    "warten",
    "irren",
    "wüten",
    "heilen",
    "pflegen",
    "betteln",
    "schmecken",
    "krachen",
    "hoffen",
    "ehren",
    "achten",
    "wenden",
    "wachsen",
    "sitzen",
    "rufen",
    "liegen",
    "wecken",
    "segeln",
    "spinnen",
    "decken",
    "graben",
    "klettern",
    "quälen",
    "äußern",
    "wagen",
    "werden",
    "werfen",
    "klingen",
    "schrecken",
    "biegen",
    "kriegen",
    "bitten",
    "brauchen",
    "stehlen",
    "wechseln",
    "lügen",
    "hören",
    "saufen",
    "fangen",
    "führen",
    "ahnen",
    "grüßen",
    "orten",
    "fesseln",
    "pfeifen",
    "jubeln",
    "wohnen",
    "lenken",
    "spielen",
    "ändern",
    "sprengen",
    "schlucken",
    "zögern",
    "leisten",
    "leugnen",
    "spüren",
    "heulen",
    "siegen",
    "schwören",
    "stammen",
    "kümmern",
    "stoppen",
    "schauen",
    "schwächen",
    "wundern",
    "schwingen",
    "machen",
    "schämen",
    "schulden",
    "fahren",
    "stürmen",
    "bellen",
    "bergen",
    "backen",
    "sichern",
    "gleichen",
    "sterben",
    "trennen",
    "zünden",
    "ruhen",
    "loben",
    "schenken",
    "kennen",
    "nennen",
    "helfen",
    "formen",
    "narren",
    "planen"
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $words[$i]
}

